$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115.. down by one.
$ws.Rows.Item(115).EntireRow.Insert()

# Populate the newly inserted row 115 with its data.
$ws.Cells.Item(115, 1).Value = 5
$ws.Cells.Item(115, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(115, 3).Value = "Maule"
$ws.Cells.Item(115, 4).Value = 44572
$ws.Cells.Item(115, 5).Value = 7
$ws.Cells.Item(115, 6).Value = 100112045
$ws.Cells.Item(115, 7).Value = "Zapallo"
$ws.Cells.Item(115, 8).Value = "Camote"
$ws.Cells.Item(115, 9).Value = "1a nueva(o)"
$ws.Cells.Item(115, 10).Value = 900
$ws.Cells.Item(115, 11).Value = 300
$ws.Cells.Item(115, 12).Value = 300
$ws.Cells.Item(115, 13).Value = 300
$ws.Cells.Item(115, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(115, 15).Value = "Región del Maule"
$ws.Cells.Item(115, 16).Value = 300
$ws.Cells.Item(115, 17).Value = 1
$ws.Cells.Item(115, 18).Value = "Hortaliza"

Write-Output "done"
